$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 6.677397333333334
$ws.Range("H2").Value = 20.032192
$ws.Range("I2").Value = 0.01905132696318048
$ws.Range("J2").Value = 0.01905132696318048
$ws.Range("M2").Value = 7.616606666666667
$ws.Range("N2").Value = 22.84982
$ws.Range("O2").Value = 0.2946797543509583
$ws.Range("P2").Value = 0.2946797543509583
$ws.Range("Q2").Value = 50.8591090450489
$ws.Range("R2").Value = 457.7319814054401
$ws.Range("S2").Value = 0.005614040349569813
$ws.Range("T2").Value = 0.005614040349569814
$ws.Range("G3").Value = 6.677397333333334
$ws.Range("H3").Value = 20.032192
$ws.Range("I3").Value = 0.01905132696318048
$ws.Range("J3").Value = 0.01905132696318048
$ws.Range("O3").Value = 0.3968845127662079
$ws.Range("P3").Value = 0.3968845127662079
$ws.Range("Q3").Value = 68.49874283873424
$ws.Range("R3").Value = 616.488685548608
$ws.Range("S3").Value = 0.007561176619331605
$ws.Range("T3").Value = 0.007561176619331605
$ws.Range("G4").Value = 6.677397333333334
$ws.Range("H4").Value = 20.032192
$ws.Range("I4").Value = 0.01905132696318048
$ws.Range("J4").Value = 0.01905132696318048
$ws.Range("O4").Value = 0.3084357328828338
$ws.Range("P4").Value = 0.3084357328828339
$ws.Range("Q4").Value = 53.23326879591112
$ws.Range("R4").Value = 479.0994191632
$ws.Range("S4").Value = 0.005876109994279065
$ws.Range("T4").Value = 0.005876109994279066
$ws.Range("I5").Value = 0.8885011423915244
$ws.Range("J5").Value = 0.8885011423915244
$ws.Range("M5").Value = 7.616606666666667
$ws.Range("N5").Value = 22.84982
$ws.Range("O5").Value = 0.2946797543509583
$ws.Range("P5").Value = 0.2946797543509583
$ws.Range("Q5").Value = 2371.92803288056
$ws.Range("R5").Value = 21347.35229592504
$ws.Range("S5").Value = 0.2618232983804802
$ws.Range("T5").Value = 0.2618232983804803
$ws.Range("I6").Value = 0.8885011423915244
$ws.Range("J6").Value = 0.8885011423915244
$ws.Range("O6").Value = 0.3968845127662079
$ws.Range("P6").Value = 0.3968845127662079
$ws.Range("S6").Value = 0.3526323429902793
$ws.Range("T6").Value = 0.3526323429902793
$ws.Range("I7").Value = 0.8885011423915244
$ws.Range("J7").Value = 0.8885011423915244
$ws.Range("O7").Value = 0.3084357328828338
$ws.Range("P7").Value = 0.3084357328828339
$ws.Range("S7").Value = 0.2740455010207649
$ws.Range("T7").Value = 0.274045501020765
$ws.Range("G8").Value = 32.40240933333334
$ws.Range("H8").Value = 97.20722800000001
$ws.Range("I8").Value = 0.09244753064529498
$ws.Range("J8").Value = 0.09244753064529498
$ws.Range("M8").Value = 7.616606666666667
$ws.Range("N8").Value = 22.84982
$ws.Range("O8").Value = 0.2946797543509583
$ws.Range("P8").Value = 0.2946797543509583
$ws.Range("Q8").Value = 246.7964069443289
$ws.Range("R8").Value = 2221.16766249896
$ws.Range("S8").Value = 0.02724241562090821
$ws.Range("T8").Value = 0.02724241562090822
$ws.Range("G9").Value = 32.40240933333334
$ws.Range("H9").Value = 97.20722800000001
$ws.Range("I9").Value = 0.09244753064529498
$ws.Range("J9").Value = 0.09244753064529498
$ws.Range("O9").Value = 0.3968845127662079
$ws.Range("P9").Value = 0.3968845127662079
$ws.Range("Q9").Value = 332.3936248633303
$ws.Range("R9").Value = 2991.542623769973
$ws.Range("S9").Value = 0.03669099315659697
$ws.Range("T9").Value = 0.03669099315659697
$ws.Range("G10").Value = 32.40240933333334
$ws.Range("H10").Value = 97.20722800000001
$ws.Range("I10").Value = 0.09244753064529498
$ws.Range("J10").Value = 0.09244753064529498
$ws.Range("O10").Value = 0.3084357328828338
$ws.Range("P10").Value = 0.3084357328828339
$ws.Range("Q10").Value = 258.3171375868112
$ws.Range("S10").Value = 0.0285141218677898
$ws.Range("T10").Value = 0.0285141218677898
